$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first occurrence of updated values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3310
$wsExpo.Range("F3").Value = 741
$wsExpo.Range("F5").Value = 6928
$wsExpo.Range("F6").Value = 2263
$wsExpo.Range("F8").Value = 89
$wsExpo.Range("F9").Value = 27
$wsExpo.Range("F12").Value = 28
$wsExpo.Range("F13").Value = 162
$wsExpo.Range("F14").Value = 288
$wsExpo.Range("F15").Value = 40

# Sheet "全部类型" (All types) - second occurrence of the same updated values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3310
$wsAll.Range("F4").Value = 741
$wsAll.Range("F6").Value = 6928
$wsAll.Range("F7").Value = 2263
$wsAll.Range("F9").Value = 89
$wsAll.Range("F10").Value = 27
$wsAll.Range("F13").Value = 28
$wsAll.Range("F14").Value = 162
$wsAll.Range("F15").Value = 288
$wsAll.Range("F16").Value = 40
